# Loan RBI, Variable Instalments
# - Insert a new (blank) column before column N on the "Repayment Schedule"
#   sheet, pushing the existing "Late" / heading / "Outstanding" columns one
#   slot to the right.
# - Switch the active/selected sheet from "NewLoanInput" to
#   "Repayment Schedule", leaving a fresh selection there.

$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")

# Insert a blank column in front of column N (shifts N->O, O->P, P->Q, ...).
$wsSchedule.Columns("N:N").Insert()

# Make "Repayment Schedule" the active sheet/tab and park the selection on
# K20, matching the saved view state.
$wsSchedule.Activate()
$wsSchedule.Range("K20").Select()
